$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new data row right above the current row 591.
# Excel shifts row 591..690 down to 592..691 automatically.
$ws.Rows("591").Insert()

# Populate the newly inserted row 591 with the new record.
# Columns that are identical to the (now shifted-down) original row stay the same,
# only the columns called out in the diff get the new values.
$ws.Range("A591").Value = 3
$ws.Range("B591").Value = 'Femacal de La Calera'
$ws.Range("C591").Value = 'Coquimbo'
$ws.Range("D591").Value = 45218
$ws.Range("E591").Value = 5
$ws.Range("F591").Value = 100112031
$ws.Range("G591").Value = 'Poroto verde'
$ws.Range("H591").Value = 'Sin especificar'
$ws.Range("I591").Value = 'Primera'
$ws.Range("J591").Value = 56
$ws.Range("K591").Value = 30000
$ws.Range("L591").Value = 30000
$ws.Range("M591").Value = 30000
$ws.Range("N591").Value = '$/malla 25 kilos'
$ws.Range("O591").Value = 'Provincia de Limarí'
$ws.Range("P591").Value = 1200
$ws.Range("Q591").Value = 25
$ws.Range("R591").Value = 'Hortaliza'

# Make sure the date cell keeps the workbook's date-time number format.
$ws.Range("D591").NumberFormat = $ws.Range("D592").NumberFormat
